$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update header / label text (shared string content change) ---
# This single text change causes the shared-string table to drop the now-unused
# "CUSTO POR HORA DO PAD" entry and append the new text at the end, which in turn
# shifts every other shared-string index used elsewhere on the sheet - matching
# the target diff automatically.
$ws.Range("A26").Value = "CUSTO DE TRABALHO DE 1 SERVIDOR POR HORA NO PAD"

# --- Fix broken formulas (previously referencing #REF!) ---
$ws.Range("D19").Formula = "=B3+D12"
$ws.Range("D20").Formula = "=D3*3"
$ws.Range("D21").Formula = "=D3+E12"
$ws.Range("E21").Formula = "=D21*30"

# D24, D25, D26 keep their original formulas; once D19/D20/D21 resolve, their
# #REF! errors clear up automatically on recalculation.

# --- Add a top border to B20:C20 to match the A20 label box style used elsewhere ---
$ws.Range("B20:C20").Borders.Item(8).LineStyle = 1
$ws.Range("B20:C20").Borders.Item(8).Weight = 2
$ws.Range("B20:C20").Borders.Item(8).ColorIndex = 1

# --- Re-order the A19:C19 / A20:C20 merged-cell bookkeeping entries ---
$ws.Range("A19:C19").UnMerge()
$ws.Range("A20:C20").UnMerge()
$ws.Range("A20:C20").Merge()
$ws.Range("A19:C19").Merge()

# --- Update the view selection ---
$ws.Activate()
$ws.Range("D9").Select()

$excel.Calculate()
